$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.329099999999999
$ws.Range("D5").Value = -7.719900000000004
$ws.Range("A10").Value = -22.00869999999999
$ws.Range("A12").Value = -21.7637
$ws.Range("B12").Value = 6.026099999999996
$ws.Range("C12").Value = -11.0678
$ws.Range("C13").Value = -12.5678
$ws.Range("B17").Value = 5.3217
$ws.Range("A18").Value = -22.15550000000001
$ws.Range("C21").Value = -13.77959999999999
$ws.Range("D23").Value = -8.235800000000003
$ws.Range("B26").Value = 4.476500000000002
$ws.Range("D26").Value = -7.874000000000001
$ws.Range("B27").Value = 5.5932
$ws.Range("B28").Value = 6.191600000000001
$ws.Range("D34").Value = -7.671700000000003
$ws.Range("C36").Value = -11.923
$ws.Range("A37").Value = -21.91249999999999
$ws.Range("B37").Value = 6.3811
$ws.Range("C38").Value = -13.0417
$ws.Range("D40").Value = -8.995299999999991
$ws.Range("C41").Value = -12.71200000000001
$ws.Range("D47").Value = -7.763199999999999
$ws.Range("C52").Value = -11.05470000000001
$ws.Range("A55").Value = -22.2032
$ws.Range("D60").Value = -8.617699999999999
$ws.Range("B65").Value = 6.316500000000003
$ws.Range("C67").Value = -11.5319
$ws.Range("A68").Value = -21.47739999999999
$ws.Range("D72").Value = -6.825499999999998
$ws.Range("B73").Value = 9.444800000000001
$ws.Range("A77").Value = -20.73299999999999
$ws.Range("A78").Value = -19.71579999999997
$ws.Range("D83").Value = -8.9902
$ws.Range("B84").Value = 5.0575
$ws.Range("B85").Value = 5.912100000000004
$ws.Range("C89").Value = -14.17779999999999
$ws.Range("B93").Value = 5.5837
$ws.Range("B95").Value = 6.549700000000005
$ws.Range("C95").Value = -12.42339999999999
$ws.Range("B98").Value = 5.345100000000006
$ws.Range("B99").Value = 5.594699999999996
$ws.Range("B101").Value = 5.7759
$ws.Range("C105").Value = -12.40930000000001
